# The workbook tracks a list of scraped URLs (and related metadata) in
# column A of the "urls" sheet. Three URLs that no longer belong in the
# list are being removed:
#   - businesswire.com Silicon-Ranch article   (was row 6)
#   - powerfinancerisk.com article             (was row 7)
#   - marathoncapital.com transactions/spower  (was row 9)
# The dbusiness.com row (was row 8) is kept and shifts up to row 6.
#
# Delete the rows bottom-up by their original row numbers so each
# Delete() call targets the intended row regardless of prior shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(9).Delete()
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()

# Mirror the author's final selection (was A9, now A6 after the deletes).
$ws.Range("A6").Select()
